$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the per-row data that gets rotated/swapped between rows:
# A, B, D, E, F, G, H, M, Q, R
# (M is sometimes absent; an absent M is represented here as an empty string "",
#  which - when written back - removes the cell's content entirely, matching the diff.)

function Get-RowData($r) {
    $row = @{}
    $row.A = $ws.Cells.Item($r, 1).Value2
    $row.B = $ws.Cells.Item($r, 2).Value2
    $row.D = $ws.Cells.Item($r, 4).Value2
    $row.E = $ws.Cells.Item($r, 5).Value2
    $row.F = $ws.Cells.Item($r, 6).Value2
    $row.G = $ws.Cells.Item($r, 7).Value2
    $row.H = $ws.Cells.Item($r, 8).Value2
    $m = $ws.Cells.Item($r, 13).Value2
    if ($null -eq $m) { $m = "" }
    $row.M = $m
    $row.Q = $ws.Cells.Item($r, 17).Value2
    $row.R = $ws.Cells.Item($r, 18).Value2
    return $row
}

function Set-RowData($r, $data) {
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
}

# Snapshot the "before" values of every affected row first, so writes don't clobber reads.
$row2  = Get-RowData 2
$row3  = Get-RowData 3
$row5  = Get-RowData 5
$row14 = Get-RowData 14
$row15 = Get-RowData 15
$row16 = Get-RowData 16
$row17 = Get-RowData 17
$row20 = Get-RowData 20
$row22 = Get-RowData 22
$row27 = Get-RowData 27
$row28 = Get-RowData 28

# 3-cycle: row2 <- old row3, row3 <- old row5, row5 <- old row2
Set-RowData 2 $row3
Set-RowData 3 $row5
Set-RowData 5 $row2

# 2-cycles (simple swaps)
Set-RowData 14 $row15
Set-RowData 15 $row14

Set-RowData 16 $row17
Set-RowData 17 $row16

Set-RowData 20 $row22
Set-RowData 22 $row20

Set-RowData 27 $row28
Set-RowData 28 $row27
